# delete and reset functionality added
# Append two new bill rows (a fresh customer bill + a repeat "fanta" order)
# to the bottom of the bills worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text looks like a date or a plain number must be forced to
# Text format before assignment, otherwise Excel auto-coerces them into a
# date serial / number instead of keeping them as the literal strings used
# throughout the rest of this sheet.
$textRefs = @("A14","D14","F14","H14","I14","A15","F15","H15","I15")
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row 14: a brand-new customer bill (name/address/contact populated)
$ws.Range("A14").Value = "2019-11-12"
$ws.Range("B14").Value = "sumant gupta"
$ws.Range("C14").Value = "matihani"
$ws.Range("D14").Value = "9824840876"
$ws.Range("E14").Value = "[sugar - 1, suji - 20]"
$ws.Range("F14").Value = "870"
$ws.Range("G14").Value = 87.0
$ws.Range("H14").Value = "101.79"
$ws.Range("I14").Value = "884.79"

# Row 15: a repeat entry (only date + particulars populated, like the other
# continuation rows already on the sheet)
$ws.Range("A15").Value = "2019-11-14"
$ws.Range("E15").Value = "[fanta - 5]"
$ws.Range("F15").Value = "950"
$ws.Range("G15").Value = 47.499999999999886
$ws.Range("H15").Value = "117.33"
$ws.Range("I15").Value = "1019.83"

# Restore default (General) formatting now that the text values are locked
# in, so the new rows don't carry any visible/explicit number format.
foreach ($ref in $textRefs) {
    $ws.Range($ref).ClearFormats()
}
